# Automatische test-sync: 2025-08-28 20:21:50
# Adds the newest "Retour status" log entry (row 13) to the Logs sheet and
# bumps the matching Dashboard aantal count from 11 to 12.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 13 -----------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A13").Value = "Retour status"
$logs.Range("B13").Value = "mailmind.test@zohomail.eu"
$logs.Range("D13").Value = "Retour / Terugbetaling"
$logs.Range("F13").Value = "2025-08-28 20:21:11"
$logs.Range("G13").Value = "Ja"
$logs.Range("H13").Value = "Nee"
$logs.Range("I13").Value = "Nee"
$logs.Range("J13").Value = "Nee"

# --- Extend the conditional-formatting ranges to include the new row ---
$logs.Range("D2:D12").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D13"))
$logs.Range("G2:G12").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G13"))
$logs.Range("H2:H12").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H13"))
$logs.Range("I2:I12").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I13"))
$logs.Range("J2:J12").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J13"))

# --- Dashboard sheet: bump the "Retour / Terugbetaling" count ----------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 12
